$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells P1, Q1 with values 14 and 15, using same style as existing header row (B1:O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Copy style from O1 (last existing header cell) to P1:Q1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# For every data row (2-25): swap I<->K and M<->O values, then add P and Q columns = 2
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Cells.Item($r, 9).Value()   # column I = 9
    $kVal = $ws.Cells.Item($r, 11).Value()  # column K = 11
    $mVal = $ws.Cells.Item($r, 13).Value()  # column M = 13
    $oVal = $ws.Cells.Item($r, 15).Value()  # column O = 15

    $ws.Cells.Item($r, 9).Value = $kVal
    $ws.Cells.Item($r, 11).Value = $iVal
    $ws.Cells.Item($r, 13).Value = $oVal
    $ws.Cells.Item($r, 15).Value = $mVal

    $ws.Cells.Item($r, 16).Value = 2  # column P = 16
    $ws.Cells.Item($r, 17).Value = 2  # column Q = 17
}
